# Fix 1620 point values in grading scale (Sheet2 "Letter Grade" table).
# The minimum-point thresholds were mis-entered ~100 pts too high for the
# bottom rows; correct them, including the "<480" label that described the
# bottom cut-off.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Value = 600
$ws2.Range("B3").Value = 590
$ws2.Range("B4").Value = 540
$ws2.Range("B5").Value = 530
$ws2.Range("B6").Value = 460
$ws2.Range("B7").Value = 450
$ws2.Range("B8").Value = 380
$ws2.Range("B9").Value = "<380"

# Page numbering tweak that accompanied the fix.
$ws2.PageSetup.FirstPageNumber = 0

# Leave the cursor where the author ended up after the edit.
$ws2.Activate()
$ws2.Range("B10").Select()
